$wb = $excel.ActiveWorkbook

# Worksheets (by their tab order in the workbook):
#   Item(1) = "1-4"  (table_details 1-4 sheet)
#   Item(2) = "5-8"  (table_details 5-8 sheet)
$ws14 = $wb.Worksheets.Item(1)
$ws58 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Sheet "1-4": the "work_place" table (rows 25-31). Its description column
# ("详细工作地点" / D) switches from nchar to nvarchar, and the fixed-length
# value that only applies to nchar is cleared out.
# ---------------------------------------------------------------------------
$ws14.Range("D28").Value = "nvarchar"
$ws14.Range("D29").ClearContents()

# ---------------------------------------------------------------------------
# Sheet "5-8": three tables get their description column's type switched
# from nchar to nvarchar, with the matching length cell cleared.
# ---------------------------------------------------------------------------
# job_type table (rows 1-7): "职能类别描述" / D4, length D5
$ws58.Range("D4").Value = "nvarchar"
$ws58.Range("D5").ClearContents()

# emergency_degree table (rows 9-15): "紧急度描述" / D12, length D13
$ws58.Range("D12").Value = "nvarchar"
$ws58.Range("D13").ClearContents()

# stuff_type table (rows 25-31): description length (D29) is cleared
$ws58.Range("D29").ClearContents()

# ---------------------------------------------------------------------------
# View state: the workbook's active tab moves from "1-4" to "5-8", and the
# selection on each sheet moves as well.
# ---------------------------------------------------------------------------
$ws14.Range("F31").Select()

$ws58.Activate()
$ws58.Range("F33").Select()
